# UC1 reviewed with Sebastian + UC1-AD added
#
# Inserts a new "Aktivitetsdiagram" (Activity Diagram) section right after
# the existing "# : 1,2,3.." / "NOT '01'" example row on the "navngivning"
# sheet, pushing every subsequent section down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("navngivning")

# Make room: insert 3 blank rows above the old row 8 ("Usecase" header).
$ws.Rows("7:9").Insert()

# New section header (row 7) - copy the yellow header style used by the
# other section headers (e.g. row 4, "Usecase").
$ws.Range("A4:B4").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)
$ws.Range("A7").Value = "Aktivitetsdiagram"

# The row-insert carried over formatted-but-empty D/E cells from row 6;
# row 7 doesn't need them.
$ws.Range("D7").Clear()
$ws.Range("E7").Clear()

# New data row (row 8) describing the AD naming convention.
$ws.Range("A8").Value = "FS-UC#-navn-AD"
$ws.Range("B8").Value = "FS-UC1-seHistorik-AD"
$ws.Range("C8").Value = "for each formel usecase"
$ws.Range("D8").Value = "in between - : no space, navn : no space, start with lowercase (verb), uppercase for (object) "
$ws.Range("E8").Clear()
$ws.Rows("8").RowHeight = 72

# Trailing example row (row 9), matching the "# : 1,2,3.." / "NOT '01'"
# pattern already used under the "Usecase" section (row 6).
$ws.Range("D9").Value = "# : 1,2,3.."
$ws.Range("E9").Value = "NOT '01'"

# Match the author's final selection/view state.
$ws.Range("B8").Select()
